$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New EPS data became available for year 1400 (reported 1401/05/22). Insert
# a fresh row right under the header so the existing history (years
# 1399..1396) shifts down one row, then fill the new row in.
$ws.Rows("2:2").Insert()
$ws.Range("A2:F2").ClearFormats()

$ws.Range("A2").Value = "'1400"
$ws.Range("B2").Value = "'1401/05/22"
$ws.Range("C2").Value = 3137
$ws.Range("D2").Value = 2150
$ws.Range("E2").Value = 700000
$ws.Range("F2").Value = 700000

$ws.Range("A2:F2").ClearFormats()
